# SimaPro - ecoinvent - biosphere mapping workbook
# Add processing for air extraction of 'Water, CN' from SimaPro exports.
#
# The "ee" sheet maps a SimaPro elementary-flow category (col A) + SimaPro
# flow name (col B) pair onto the matching ecoinvent biosphere flow name (col C).
# The "Resources" category (rows ~278-342) already enumerates one row per
# "Water, <country code>" SimaPro flow, each mapped to the generic ecoinvent
# flow "Water, unspecified natural origin" -- except for "Water, CN", whose slot
# is instead used to map the SimaPro cooling-water flow "Water, cooling,
# unspecified natural origin/m3" onto ecoinvent's "Water, cooling, unspecified
# natural origin" (the country-coded flow itself is still matched by an earlier,
# pre-existing row).
#
# This change duplicates that same 65-row block under the "Emissions to air"
# category, so that SimaPro air-emission exports of these "Water, <cc>" flows
# (including the "Water, CN" / cooling-water special case) get resolved too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ee")

# A previously blank/placeholder SimaPro category marker cell now carries the
# literal flow-name fragment used elsewhere as a lookup key.
$ws.Range("B8").Value = 'x "Wat'

# (SimaPro flow name, ecoinvent biosphere flow name) pairs to append under the
# "Emissions to air" category, mirroring the existing "Resources" block.
$newRows = @(
    @('Water, AT', 'Water, unspecified natural origin'),
    @('Water, AU', 'Water, unspecified natural origin'),
    @('Water, BA', 'Water, unspecified natural origin'),
    @('Water, BE', 'Water, unspecified natural origin'),
    @('Water, BG', 'Water, unspecified natural origin'),
    @('Water, BR', 'Water, unspecified natural origin'),
    @('Water, CA', 'Water, unspecified natural origin'),
    @('Water, CH', 'Water, unspecified natural origin'),
    @('Water, CL', 'Water, unspecified natural origin'),
    @('Water, CN', 'Water, unspecified natural origin'),
    @('Water, cooling, unspecified natural origin/m3', 'Water, cooling, unspecified natural origin'),
    @('Water, CZ', 'Water, unspecified natural origin'),
    @('Water, DE', 'Water, unspecified natural origin'),
    @('Water, DK', 'Water, unspecified natural origin'),
    @('Water, ES', 'Water, unspecified natural origin'),
    @('Water, Europe without Switzerland', 'Water, unspecified natural origin'),
    @('Water, FI', 'Water, unspecified natural origin'),
    @('Water, FR', 'Water, unspecified natural origin'),
    @('Water, GB', 'Water, unspecified natural origin'),
    @('Water, GLO', 'Water, unspecified natural origin'),
    @('Water, GR', 'Water, unspecified natural origin'),
    @('Water, HR', 'Water, unspecified natural origin'),
    @('Water, HU', 'Water, unspecified natural origin'),
    @('Water, ID', 'Water, unspecified natural origin'),
    @('Water, IE', 'Water, unspecified natural origin'),
    @('Water, IN', 'Water, unspecified natural origin'),
    @('Water, IR', 'Water, unspecified natural origin'),
    @('Water, IT', 'Water, unspecified natural origin'),
    @('Water, JP', 'Water, unspecified natural origin'),
    @('Water, KR', 'Water, unspecified natural origin'),
    @('Water, LU', 'Water, unspecified natural origin'),
    @('Water, MA', 'Water, unspecified natural origin'),
    @('Water, MX', 'Water, unspecified natural origin'),
    @('Water, MY', 'Water, unspecified natural origin'),
    @('Water, NL', 'Water, unspecified natural origin'),
    @('Water, NO', 'Water, unspecified natural origin'),
    @('Water, NORDEL', 'Water, unspecified natural origin'),
    @('Water, PE', 'Water, unspecified natural origin'),
    @('Water, PG', 'Water, unspecified natural origin'),
    @('Water, PH', 'Water, unspecified natural origin'),
    @('Water, PL', 'Water, unspecified natural origin'),
    @('Water, PT', 'Water, unspecified natural origin'),
    @('Water, RAF', 'Water, unspecified natural origin'),
    @('Water, RAS', 'Water, unspecified natural origin'),
    @('Water, RER', 'Water, unspecified natural origin'),
    @('Water, RLA', 'Water, unspecified natural origin'),
    @('Water, RME', 'Water, unspecified natural origin'),
    @('Water, RNA', 'Water, unspecified natural origin'),
    @('Water, RO', 'Water, unspecified natural origin'),
    @('Water, RoW', 'Water, unspecified natural origin'),
    @('Water, RS', 'Water, unspecified natural origin'),
    @('Water, RU', 'Water, unspecified natural origin'),
    @('Water, SE', 'Water, unspecified natural origin'),
    @('Water, SI', 'Water, unspecified natural origin'),
    @('Water, SK', 'Water, unspecified natural origin'),
    @('Water, TH', 'Water, unspecified natural origin'),
    @('Water, TR', 'Water, unspecified natural origin'),
    @('Water, TW', 'Water, unspecified natural origin'),
    @('Water, TZ', 'Water, unspecified natural origin'),
    @('Water, UA', 'Water, unspecified natural origin'),
    @('Water, UCTE', 'Water, unspecified natural origin'),
    @('Water, UCTE without Germany', 'Water, unspecified natural origin'),
    @('Water, US', 'Water, unspecified natural origin'),
    @('Water, WEU', 'Water, unspecified natural origin'),
    @('Water, ZA', 'Water, unspecified natural origin')
)

$category = "Emissions to air"
$startRow = $ws.UsedRange.Rows.Count + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $category
    $ws.Cells.Item($r, 2).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][1]
}

# Leave the selection where the editor ended up after appending the rows.
$ws.Activate()
$ws.Cells.Item($startRow + $newRows.Count + 3, 1).Select()
